$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.319.85'
$ws.Range("E2").Value = '  -0.46%  '

$ws.Range("D3").Value = '1.589.49'
$ws.Range("E3").Value = '  -0.15%  '

$ws.Range("E4").Value = '  -0.70%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '

$ws.Range("E7").Value = '  -0.68%  '

$ws.Range("E8").Value = '  +0.63%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0611'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.42'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.66%  '

$ws.Range("D12").Value = '1.812.03'
$ws.Range("E12").Value = '  -0.19%  '

$ws.Range("D13").Value = '1.592.98'
$ws.Range("E13").Value = '  +0.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.21%  '

$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.54%  '

$ws.Range("D17").Value = '26.321.56'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("D18").Value = '0.0₃0730'
$ws.Range("E18").Value = '  -0.71%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.50'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.61%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.20%  '

$ws.Range("E21").Value = '  -0.64%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.96'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '

$ws.Range("E26").Value = '  -0.58%  '

$ws.Range("E27").Value = '  +0.34%  '

$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.27%  '

$ws.Range("E30").Value = '  +0.73%  '

$ws.Range("E31").Value = '  +0.26%  '

$ws.Range("E32").Value = '  -0.61%  '

$ws.Range("E33").Value = '  +2.18%  '

$ws.Range("D34").Value = '1.319.56'
$ws.Range("E34").Value = '  +2.96%  '

$ws.Range("E35").Value = '  -1.83%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.609'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.21%  '

$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("E39").Value = '  -11.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("E41").Value = '  -0.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.76%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.766'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.45%  '

$ws.Range("E44").Value = '  -0.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("D46").Value = '1.724.77'
$ws.Range("E46").Value = '  -0.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '87.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.49'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.74%  '

$ws.Range("E49").Value = '  -1.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0978'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.95%  '

$ws.Range("E51").Value = '  -0.63%  '
